# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" column (G) values for each game row (2-38)
$kValues = @{
    2  = 2
    3  = 1
    4  = 4
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 3
    30 = 3
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 3
    36 = 1
    37 = 2
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
